$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the two trailing rows (25 and 26) entirely -------------------
# Delete bottom-up so row numbers above aren't disturbed mid-way.
$ws.Range("A26").EntireRow.Delete()
$ws.Range("A25").EntireRow.Delete()

# --- 2. Clear cells that must become empty ----------------------------------
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()
$ws.Range("A24").Clear()

# --- 3. Write the new / relocated cell values --------------------------------
$ws.Range("B10").Value = "5840897 - Clodoaldo Saron"
$ws.Range("C10").Value = "5840897 - Clodoaldo Saron"

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "15/07/2016"
$ws.Range("C13").Value = "15/07/2016"

$ws.Range("A14").Value = "Short syllabus:"

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "5840897 - Clodoaldo Saron"
$ws.Range("C15").Value = "5840897 - Clodoaldo Saron"

$ws.Range("A16").Value = "Syllabus:"

$ws.Range("A17").Value = "Avaliação:"

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C18").Value = "1033242 - Fábio Herbst Florenzano"

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento."
$ws.Range("C19").Value = "Experimentos desenvolvidos em laboratório didático; realização de relatórios para cada experimento."

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média aritmética das notas obtidas nos relatórios. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."
$ws.Range("C20").Value = "Média aritmética das notas obtidas nos relatórios. Será aprovado o aluno que obtiver nota final maior ou igual a 5,0."

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Devido às características práticas da disciplina, não será oferecida recuperação."
$ws.Range("C21").Value = "Devido às características práticas da disciplina, não será oferecida recuperação."

$ws.Range("A22").Value = "Requisitos:"

$ws.Range("B23").Value = "LOM3057 -  Introdução aos Materiais Poliméricos  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOM3057 -  Introdução aos Materiais Poliméricos  (Requisito fraco)`n"

$ws.Range("B24").Value = "LOM3058 -  Química de Polímeros  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOM3058 -  Química de Polímeros  (Requisito fraco)`n"

# --- 4. Fix up row heights ----------------------------------------------------
$ws.Rows(13).RowHeight = 60
$ws.Rows(14).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(16).RowHeight = 120
$ws.Rows(17).AutoFit()
$ws.Rows(18).RowHeight = 60
$ws.Rows(19).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(22).AutoFit()
$ws.Rows(23).RowHeight = 30
$ws.Rows(24).RowHeight = 30
